$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.221.73'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.661.07'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.007'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2637'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06321'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.66'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.502'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.667.26'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '1.889.03'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5560'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '0.0₅8024'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.31'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').Value = '26.236.33'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.653'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '197.00'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.972'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.008'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.01'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1209'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.169'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.514'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05784'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.73%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.279'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.491'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.352'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.586'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.75%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.810'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9547'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.423'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5747'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01593'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.978'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.88%  '
$ws.Range('D41').Value = '1.063.17'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8534'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.007'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '103.12'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('D45').Value = '1.799.82'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '58.32'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.011'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4411'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₈103'
$ws.Range('E49').Value = '  -2.73%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.003'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05197'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.68%  '
